$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $orig = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $orig
}

Set-TextValue $ws.Range("D2") "30.306.20"
$ws.Range("E2").Value = "  +0.04%  "

Set-TextValue $ws.Range("D3") "1.931.52"
$ws.Range("E3").Value = "  +0.15%  "

Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.05%  "

Set-TextValue $ws.Range("D5") "0.7494"
$ws.Range("E5").Value = "  +4.53%  "

Set-TextValue $ws.Range("D6") "243.05"
$ws.Range("E6").Value = "  -2.38%  "

$ws.Range("E7").Value = "  +0.05%  "

Set-TextValue $ws.Range("D8") "27.65"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("E9").Value = "  -1.23%  "

Set-TextValue $ws.Range("D10") "0.07117"
$ws.Range("E10").Value = "  +0.23%  "

Set-TextValue $ws.Range("D11") "0.08059"
$ws.Range("E11").Value = "  +0.66%  "

Set-TextValue $ws.Range("D12") "0.7796"
$ws.Range("E12").Value = "  -1.46%  "

Set-TextValue $ws.Range("D13") "1.909.45"
$ws.Range("E13").Value = "  -1.07%  "

Set-TextValue $ws.Range("D14") "5.399"
$ws.Range("E14").Value = "  +0.13%  "

Set-TextValue $ws.Range("D15") "93.06"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("E16").Value = "  -0.73%  "

Set-TextValue $ws.Range("D17") "30.308.55"
$ws.Range("E17").Value = "  +0.00%  "

Set-TextValue $ws.Range("D18") "6.021"
$ws.Range("E18").Value = "  +4.28%  "

Set-TextValue $ws.Range("D19") "251.70"
$ws.Range("E19").Value = "  -2.34%  "

Set-TextValue $ws.Range("D20") "0.000007937"
$ws.Range("E20").Value = "  -1.60%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D21") "1.001"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D22") "2.160.42"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("E23").Value = "  +0.04%  "

Set-TextValue $ws.Range("D24") "6.683"
$ws.Range("E24").Value = "  -2.21%  "

Set-TextValue $ws.Range("D25") "9.554"
$ws.Range("E25").Value = "  +0.04%  "

Set-TextValue $ws.Range("D26") "165.25"
$ws.Range("E26").Value = "  +0.40%  "

Set-TextValue $ws.Range("D27") "19.08"
$ws.Range("E27").Value = "  -0.15%  "

Set-TextValue $ws.Range("D28") "0.1297"
$ws.Range("E28").Value = "  +2.41%  "

Set-TextValue $ws.Range("D29") "2.183"
$ws.Range("E29").Value = "  -3.90%  "

Set-TextValue $ws.Range("D30") "1.368"
$ws.Range("E30").Value = "  +1.08%  "

Set-TextValue $ws.Range("D31") "1.559"
$ws.Range("E31").Value = "  +2.21%  "

Set-TextValue $ws.Range("D32") "4.416"
$ws.Range("E32").Value = "  +0.37%  "

Set-TextValue $ws.Range("D33") "4.145"
$ws.Range("E33").Value = "  +0.04%  "

Set-TextValue $ws.Range("D34") "0.05231"
$ws.Range("E34").Value = "  +1.63%  "

$ws.Range("E35").Value = "  +4.13%  "

Set-TextValue $ws.Range("D36") "0.7566"
$ws.Range("E36").Value = "  +1.60%  "

Set-TextValue $ws.Range("D37") "2.780"
$ws.Range("E37").Value = "  +0.51%  "

Set-TextValue $ws.Range("D38") "0.01948"
$ws.Range("E38").Value = "  -0.79%  "

$ws.Range("E39").Value = "  +0.10%  "

Set-TextValue $ws.Range("D40") "6.510"
$ws.Range("E40").Value = "  +2.23%  "

Set-TextValue $ws.Range("D41") "78.19"
$ws.Range("E41").Value = "  -0.21%  "

Set-TextValue $ws.Range("D42") "0.4525"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("E43").Value = "  -1.15%  "

Set-TextValue $ws.Range("D44") "0.8414"
$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("E45").Value = "  +0.08%  "

Set-TextValue $ws.Range("D46") "7.694"
$ws.Range("E46").Value = "  +3.36%  "

Set-TextValue $ws.Range("D47") "9.967"
$ws.Range("E47").Value = "  +1.79%  "

Set-TextValue $ws.Range("D48") "101.55"
$ws.Range("E48").Value = "  +1.00%  "

Set-TextValue $ws.Range("D49") "37.93"
$ws.Range("E49").Value = "  +3.29%  "

Set-TextValue $ws.Range("D50") "0.1232"
$ws.Range("E50").Value = "  +7.37%  "

Set-TextValue $ws.Range("D51") "961.48"
$ws.Range("E51").Value = "  +2.15%  "
